$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null
$ws.Range("H48").Value = 8900
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 8900
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 26700
$ws.Range("M48").Value = $null
$ws.Range("N48").Value = -27284
$ws.Range("H51").Value = 3293.2856
$ws.Range("I51").Value = 2867.3333
$ws.Range("J51").Value = 3409.4546
$ws.Range("K51").Value = 2867.3333
$ws.Range("L51").Value = 3409.4546
$ws.Range("M51").Value = -2383.3333
$ws.Range("N51").Value = -4377.4546
$ws.Range("H56").Value = 8900
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 8900
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 26700
$ws.Range("M56").Value = $null
$ws.Range("N56").Value = -27768
$ws.Range("H58").Value = 2350
$ws.Range("I58").Value = 700
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 2100
$ws.Range("L58").Value = 12000
$ws.Range("M58").Value = -1950
$ws.Range("N58").Value = -12300
$ws.Range("H64").Value = 25002542
$ws.Range("I64").Value = 34484788
$ws.Range("J64").Value = 3896.3635
$ws.Range("K64").Value = 34484788
$ws.Range("L64").Value = 3896.3635
$ws.Range("M64").Value = -34484540
$ws.Range("N64").Value = -4392.363499999999
$ws.Range("H67").Value = 25002542
$ws.Range("I67").Value = 34484788
$ws.Range("J67").Value = 3896.3635
$ws.Range("K67").Value = 34484788
$ws.Range("L67").Value = 3896.3635
$ws.Range("M67").Value = -34483930
$ws.Range("N67").Value = -5612.363499999999
$ws.Range("H70").Value = 5160
$ws.Range("I70").Value = 5692
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 17076
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -16806
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 5160
$ws.Range("I73").Value = 5692
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 17076
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -16140
$ws.Range("N73").Value = -9372
$ws.Range("H76").Value = 4702.143
$ws.Range("I76").Value = 2976.4707
$ws.Range("J76").Value = 6331.9443
$ws.Range("K76").Value = 2976.4707
$ws.Range("L76").Value = 6331.9443
$ws.Range("M76").Value = -2661.4707
$ws.Range("N76").Value = -6961.9443
$ws.Range("H79").Value = 4702.143
$ws.Range("I79").Value = 2976.4707
$ws.Range("J79").Value = 6331.9443
$ws.Range("K79").Value = 2976.4707
$ws.Range("L79").Value = 6331.9443
$ws.Range("M79").Value = -1884.4707
$ws.Range("N79").Value = -8515.944299999999
$ws.Range("H82").Value = 1446.25
$ws.Range("I82").Value = 279.33334
$ws.Range("J82").Value = 4947
$ws.Range("K82").Value = 838.0000200000001
$ws.Range("L82").Value = 14841
$ws.Range("M82").Value = -432.0000200000001
$ws.Range("N82").Value = -15653
$ws.Range("H85").Value = 1446.25
$ws.Range("I85").Value = 279.33334
$ws.Range("J85").Value = 4947
$ws.Range("K85").Value = 838.0000200000001
$ws.Range("L85").Value = 14841
$ws.Range("M85").Value = 565.9999799999999
$ws.Range("N85").Value = -17649
$ws.Range("H100").Value = 2929.2856
$ws.Range("I100").Value = 1851.5555
$ws.Range("J100").Value = 4869.2
$ws.Range("K100").Value = 1851.5555
$ws.Range("L100").Value = 4869.2
$ws.Range("M100").Value = -1310.5555
$ws.Range("N100").Value = -5951.2
$ws.Range("H103").Value = 521.7273
$ws.Range("I103").Value = 509.8889
$ws.Range("J103").Value = 575
$ws.Range("K103").Value = 1529.6667
$ws.Range("L103").Value = 1725
$ws.Range("M103").Value = -943.6667
$ws.Range("N103").Value = -2897
$ws.Range("H106").Value = 1800.625
$ws.Range("I106").Value = 1681
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 1681
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -1050
$ws.Range("N106").Value = -3262
$ws.Range("H118").Value = 909.225
$ws.Range("I118").Value = 589.0909
$ws.Range("J118").Value = 1030.6552
$ws.Range("K118").Value = 1767.2727
$ws.Range("L118").Value = 3091.9656
$ws.Range("M118").Value = -110.2727
$ws.Range("N118").Value = -6405.9656
$ws.Range("H129").Value = 879.375
$ws.Range("I129").Value = 729.38464
$ws.Range("J129").Value = 982
$ws.Range("K129").Value = 2188.15392
$ws.Range("L129").Value = 2946
$ws.Range("M129").Value = 2811.84608
$ws.Range("N129").Value = -12946

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 9633.333000000001
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 9450
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 9450
$ws.Range("M44").Value = -9512
$ws.Range("N44").Value = -10426
$ws.Range("H55").Value = 20000
$ws.Range("I55").Value = 20000
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -19685
$ws.Range("N55").Value = -20630
$ws.Range("H63").Value = 1720
$ws.Range("I63").Value = 1720
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1720
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1034
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 1720
$ws.Range("I66").Value = 1720
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8600
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5168
$ws.Range("N66").Value = $null
$ws.Range("H69").Value = 80000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 80000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 80000
$ws.Range("N69").Value = -81498
$ws.Range("H72").Value = 80000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 80000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 240000
$ws.Range("N72").Value = -247488
$ws.Range("H80").Value = 43500
$ws.Range("I80").Value = 24000
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 24000
$ws.Range("L80").Value = 50000
$ws.Range("M80").Value = -23002
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 43500
$ws.Range("I83").Value = 24000
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 72000
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -67008
$ws.Range("N83").Value = -159984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 48000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 48000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -52492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 974.9167
$ws.Range("I24").Value = 649.5
$ws.Range("J24").Value = 1040
$ws.Range("K24").Value = 1948.5
$ws.Range("L24").Value = 3120
$ws.Range("M24").Value = -1718.5
$ws.Range("N24").Value = -3580
$ws.Range("H32").Value = 207757.36
$ws.Range("I32").Value = 322233.66
$ws.Range("J32").Value = 1700
$ws.Range("K32").Value = 966700.98
$ws.Range("L32").Value = 5100
$ws.Range("M32").Value = -966417.98
$ws.Range("N32").Value = -5666
$ws.Range("H34").Value = 564.3461
$ws.Range("I34").Value = 272.75
$ws.Range("J34").Value = 814.2857
$ws.Range("K34").Value = 818.25
$ws.Range("L34").Value = 2442.8571
$ws.Range("M34").Value = -734.25
$ws.Range("N34").Value = -2610.8571
$ws.Range("H39").Value = 3200
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3200
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 9600
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -10188
$ws.Range("H55").Value = 3000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 727.8889
$ws.Range("I97").Value = 727.8889
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 727.8889
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -231.8889
$ws.Range("N97").Value = $null
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2905
$ws.Range("I61").Value = 2135
$ws.Range("J61").Value = 3097.5
$ws.Range("K61").Value = 2135
$ws.Range("L61").Value = 3097.5
$ws.Range("M61").Value = -1933
$ws.Range("N61").Value = -3501.5
$ws.Range("H113").Value = 2905
$ws.Range("I113").Value = 2135
$ws.Range("J113").Value = 3097.5
$ws.Range("K113").Value = 2135
$ws.Range("L113").Value = 3097.5
$ws.Range("M113").Value = 35
$ws.Range("N113").Value = -7437.5
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 47900
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 47900
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 47900
$ws.Range("N116").Value = -57078
